$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sitzungsplan")

# Insert a new row for class "06F" right after the "06E" row (currently row 13),
# copying Sitzungsleiter/Von/Bis/Datum from that row.
$ws.Rows.Item(14).Insert()
$ws.Range("A14").Value = "06F"
$ws.Range("B14").Value = "Livia Schleßing, OStRin"
$ws.Range("C14").Value = "15:00"
$ws.Range("D14").Value = "15:30"
$ws.Range("E14").Value = "18.07.2017"

# Insert a new row for class "08D" right after the "08C" row (now row 23),
# copying Sitzungsleiter/Von/Bis/Datum from that row.
$ws.Rows.Item(24).Insert()
$ws.Range("A24").Value = "08D"
$ws.Range("B24").Value = "Ute Badum, OStRin"
$ws.Range("C24").Value = "17:30"
$ws.Range("D24").Value = "18:00"
$ws.Range("E24").Value = "18.07.2017"

# Update view state: selection and scroll position.
$ws.Range("A25").Select()
$ws.Application.ActiveWindow.ScrollRow = 10
